$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the last column header from "seed_weight" to "percent_moisture"
$ws.Range("O1").Value = "percent_moisture"

# Populate the new percent_moisture column for every data row with its
# (placeholder) value
$ws.Range("O2:O37").Value = 0.00001

# Match the author's final selection/cursor position
[void]$ws.Range("R22").Select()
